# Scheduled market-data refresh: update cached Universalis price/profit
# figures (currentAveragePrice*, Leve price/profit cols H-N) across all
# eight crafting-job sheets with freshly fetched values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 748.39685
$ws.Range("J17").Value = 639.07275
$ws.Range("L17").Value = 1917.21825
$ws.Range("N17").Value = -2253.21825
$ws.Range("H98").Value = 4998.25
$ws.Range("I98").Value = 1897.375
$ws.Range("J98").Value = 11200
$ws.Range("K98").Value = 1897.375
$ws.Range("L98").Value = 11200
$ws.Range("M98").Value = -399.375
$ws.Range("N98").Value = -14196
$ws.Range("H100").Value = 13335166
$ws.Range("I100").Value = 14287321
$ws.Range("K100").Value = 14287321
$ws.Range("M100").Value = -14286780
$ws.Range("H113").Value = 3969.7646
$ws.Range("J113").Value = 4500.0713
$ws.Range("L113").Value = 4500.0713
$ws.Range("N113").Value = -11008.0713
$ws.Range("H122").Value = 4998.25
$ws.Range("I122").Value = 1897.375
$ws.Range("J122").Value = 11200
$ws.Range("K122").Value = 5692.125
$ws.Range("L122").Value = 33600
$ws.Range("M122").Value = -3242.125
$ws.Range("N122").Value = -38500
$ws.Range("H138").Value = 4944.897
$ws.Range("I138").Value = 919.36365
$ws.Range("J138").Value = 6125.72
$ws.Range("K138").Value = 2758.09095
$ws.Range("L138").Value = 18377.16
$ws.Range("M138").Value = 2381.90905
$ws.Range("N138").Value = -28657.16

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2134
$ws.Range("I122").Value = 1250
$ws.Range("K122").Value = 3750
$ws.Range("M122").Value = -1300
$ws.Range("H137").Value = 45190
$ws.Range("J137").Value = 45190
$ws.Range("L137").Value = 45190
$ws.Range("N137").Value = -55390
$ws.Range("H139").Value = 43835.7
$ws.Range("J139").Value = 43835.7
$ws.Range("L139").Value = 43835.7
$ws.Range("N139").Value = -54115.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 4565161.5
$ws.Range("J7").Value = 24086.555
$ws.Range("L7").Value = 24086.555
$ws.Range("N7").Value = -24312.555
$ws.Range("H132").Value = 50735
$ws.Range("J132").Value = 50735
$ws.Range("L132").Value = 50735
$ws.Range("N132").Value = -60855
$ws.Range("H134").Value = 4756.7437
$ws.Range("I134").Value = 1524.0714
$ws.Range("J134").Value = 12985.363
$ws.Range("K134").Value = 4572.2142
$ws.Range("L134").Value = 38956.089
$ws.Range("M134").Value = -2037.2142
$ws.Range("N134").Value = -44026.089
$ws.Range("H138").Value = 41311.92
$ws.Range("J138").Value = 41311.92
$ws.Range("L138").Value = 41311.92
$ws.Range("N138").Value = -51591.92

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3194.96
$ws.Range("I31").Value = 1239.2858
$ws.Range("J31").Value = 5684
$ws.Range("K31").Value = 1239.2858
$ws.Range("L31").Value = 5684
$ws.Range("M31").Value = -944.2858000000001
$ws.Range("N31").Value = -6274
$ws.Range("H34").Value = 3194.96
$ws.Range("I34").Value = 1239.2858
$ws.Range("J34").Value = 5684
$ws.Range("K34").Value = 1239.2858
$ws.Range("L34").Value = 5684
$ws.Range("M34").Value = -1037.2858
$ws.Range("N34").Value = -6088
$ws.Range("H58").Value = 2455.1604
$ws.Range("I58").Value = 1527.4854
$ws.Range("K58").Value = 1527.4854
$ws.Range("M58").Value = -1324.4854
$ws.Range("H86").Value = 3188
$ws.Range("I86").Value = 2615.5
$ws.Range("J86").Value = 4333
$ws.Range("K86").Value = 2615.5
$ws.Range("L86").Value = 4333
$ws.Range("M86").Value = -1492.5
$ws.Range("N86").Value = -6579
$ws.Range("H89").Value = 3188
$ws.Range("I89").Value = 2615.5
$ws.Range("J89").Value = 4333
$ws.Range("K89").Value = 13077.5
$ws.Range("L89").Value = 21665
$ws.Range("M89").Value = -7461.5
$ws.Range("N89").Value = -32897
$ws.Range("H136").Value = 2455.1604
$ws.Range("I136").Value = 1527.4854
$ws.Range("K136").Value = 4582.456200000001
$ws.Range("M136").Value = -2032.456200000001
$ws.Range("H138").Value = 43477.777
$ws.Range("J138").Value = 43477.777
$ws.Range("L138").Value = 43477.777
$ws.Range("N138").Value = -53757.777
$ws.Range("H140").Value = 87631.664
$ws.Range("J140").Value = 87631.664
$ws.Range("L140").Value = 87631.664
$ws.Range("N140").Value = -97991.664
$ws.Range("H141").Value = 30450
$ws.Range("J141").Value = 30450
$ws.Range("L141").Value = 30450
$ws.Range("N141").Value = -40810

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3095.5134
$ws.Range("I122").Value = 1033.1
$ws.Range("J122").Value = 3859.3704
$ws.Range("K122").Value = 9297.9
$ws.Range("L122").Value = 34734.3336
$ws.Range("M122").Value = -6847.9
$ws.Range("N122").Value = -39634.3336
$ws.Range("H129").Value = 2493.2632
$ws.Range("J129").Value = 2959
$ws.Range("L129").Value = 8877
$ws.Range("N129").Value = -18877
$ws.Range("H131").Value = 780.7895
$ws.Range("I131").Value = 398
$ws.Range("J131").Value = 802.05554
$ws.Range("K131").Value = 1194
$ws.Range("L131").Value = 2406.16662
$ws.Range("M131").Value = 3846
$ws.Range("N131").Value = -12486.16662
$ws.Range("H133").Value = 2639.1538
$ws.Range("I133").Value = 2459.0833
$ws.Range("J133").Value = 4800
$ws.Range("K133").Value = 7377.249899999999
$ws.Range("L133").Value = 14400
$ws.Range("M133").Value = -2317.249899999999
$ws.Range("N133").Value = -24520
$ws.Range("H138").Value = 3158
$ws.Range("J138").Value = 3150
$ws.Range("L138").Value = 9450
$ws.Range("N138").Value = -19730

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 30788.889
$ws.Range("J88").Value = 30788.889
$ws.Range("L88").Value = 30788.889
$ws.Range("N88").Value = -31690.889
$ws.Range("H91").Value = 30788.889
$ws.Range("J91").Value = 30788.889
$ws.Range("L91").Value = 30788.889
$ws.Range("N91").Value = -33908.889
$ws.Range("H126").Value = 3145.36
$ws.Range("I126").Value = 2943.2
$ws.Range("J126").Value = 3954
$ws.Range("K126").Value = 8829.599999999999
$ws.Range("L126").Value = 11862
$ws.Range("M126").Value = -6359.599999999999
$ws.Range("N126").Value = -16802
$ws.Range("H140").Value = 36908.965
$ws.Range("J140").Value = 36908.965
$ws.Range("L140").Value = 36908.965
$ws.Range("N140").Value = -47268.965

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3953.8572
$ws.Range("I136").Value = 1981.3846
$ws.Range("J136").Value = 5663.3335
$ws.Range("K136").Value = 5944.1538
$ws.Range("L136").Value = 16990.0005
$ws.Range("M136").Value = -3394.1538
$ws.Range("N136").Value = -22090.0005
$ws.Range("H139").Value = 47897.5
$ws.Range("J139").Value = 47897.5
$ws.Range("L139").Value = 47897.5
$ws.Range("N139").Value = -58177.5
$ws.Range("H140").Value = 83860.5
$ws.Range("J140").Value = 83860.5
$ws.Range("L140").Value = 83860.5
$ws.Range("N140").Value = -94220.5
$ws.Range("H141").Value = 39303
$ws.Range("J141").Value = 39303
$ws.Range("L141").Value = 39303
$ws.Range("N141").Value = -49663

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2650250.5
$ws.Range("I96").Value = 2650250.5
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 2650250.5
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -2648877.5
$ws.Range("N96").ClearContents()
$ws.Range("H122").Value = 2703.6785
$ws.Range("I122").Value = 1691.3914
$ws.Range("K122").Value = 5074.174199999999
$ws.Range("M122").Value = -2624.174199999999
